# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral"      (used by the slide master / all slides)
#   ppt/theme/theme2.xml -> "Office Theme"  (used by the notes master)
#
# The authored edit swaps their contents: theme1.xml becomes the
# "Office Theme" palette and theme2.xml becomes the "Integral" palette
# (file names / relationships are untouched, only the theme payloads
# trade places). The only part of that payload the PowerPoint object
# model can reach from slide-level automation is the *active* theme's
# color scheme (the one behind the slide master / all the slides,
# i.e. theme1.xml) via Theme.ThemeColorScheme, so re-point those 12
# slots at the Office palette.

function Convert-HexToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office theme palette, in PpColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officePalette.Count; $i++) {
    $colors.Item($i + 1).RGB = Convert-HexToBgr $officePalette[$i]
}

# Best-effort: also try to rename the theme/design/color-scheme from
# "Integral" to "Office Theme" / "Office" to match the swapped part's
# <a:theme name="..."> / <a:clrScheme name="..."> attributes.
try { $theme.Name = "Office Theme" } catch {}
try { $colors.Name = "Office" } catch {}
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
